# Weekly update: insert a new Apio (Femacal de La Calera) price record as
# row 205, pushing the existing rows 205:227 down to 206:228.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 205 (shifts 205:227 -> 206:228, carries D-column
# date style down automatically, same as Excel's native "Insert Sheet Rows").
$ws.Rows(205).Insert()

# Populate the new row with the week's data.
$ws.Cells.Item(205, 1).Value  = 3
$ws.Cells.Item(205, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(205, 3).Value  = "Coquimbo"
$ws.Cells.Item(205, 4).Value  = 44449
$ws.Cells.Item(205, 5).Value  = 5
$ws.Cells.Item(205, 6).Value  = 100112017
$ws.Cells.Item(205, 7).Value  = "Apio"
$ws.Cells.Item(205, 8).Value  = "Americana (o)"
$ws.Cells.Item(205, 9).Value  = "Primera"
$ws.Cells.Item(205, 10).Value = 250
$ws.Cells.Item(205, 11).Value = 9000
$ws.Cells.Item(205, 12).Value = 9500
$ws.Cells.Item(205, 13).Value = 9240
$ws.Cells.Item(205, 14).Value = "`$/docena de matas"
$ws.Cells.Item(205, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(205, 16).Value = 1540
$ws.Cells.Item(205, 17).Value = 6
$ws.Cells.Item(205, 18).Value = "Hortaliza"
